# Apply the metrics update edit to the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared values (V1) to apply to columns B:Q for rows 2..26
$newValues = @(
    0.445980542904741,
    -0.05067851974480586,
    -0.4794294150153033,
    0.5954233101814106,
    0.6074911698636565,
    0.3288898821261522,
    0.6237281562693088,
    0.2528224391702676,
    0.2487082662774325,
    0.25076535272385,
    0.279043433096329,
    0.5734892170966707,
    0.05025235926527027,
    0.5979038403548141,
    22.22406457805329,
    34.41282282673529
)

# Relabel the three rows that moved (A column), matching the permutation in the diff
$ws.Range("A2").Value = "model_4_7_0"
$ws.Range("A14").Value = "model_4_7_12"
$ws.Range("A26").Value = "model_4_7_24"

# Update columns B:Q for every data row (2 through 26) with the new shared values
for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
